$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.73
$ws.Range("H2").Value = 3.3
$ws.Range("I2").Value = 5.75
$ws.Range("J2").Value = 2.4
$ws.Range("K2").Value = 2.05
$ws.Range("L2").Value = 6
$ws.Range("Q2").Value = 2.35
$ws.Range("R2").Value = 1.57
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 2.5
$ws.Range("X2").Value = 7
$ws.Range("Z2").Value = 13
$ws.Range("AA2").Value = 17
$ws.Range("AE2").Value = 21
$ws.Range("AF2").Value = 81
$ws.Range("AG2").Value = 11
$ws.Range("AH2").Value = 26
$ws.Range("AI2").Value = 19
$ws.Range("AK2").Value = 51
$ws.Range("AN2").Value = 3.5
$ws.Range("AO2").Value = 9.5
$ws.Range("AP2").Value = 23
$ws.Range("AQ2").Value = 34
$ws.Range("AS2").Value = 201
$ws.Range("AT2").Value = 2.5
$ws.Range("AV2").Value = 81
$ws.Range("AW2").Value = 6.5
$ws.Range("AZ2").Value = 126

# Row 5
$ws.Range("G5").Value = 7.5
$ws.Range("H5").Value = 4.1
$ws.Range("I5").Value = 1.38
$ws.Range("J5").Value = 7.1
$ws.Range("K5").Value = 2.27
$ws.Range("L5").Value = 1.9
$ws.Range("P5").Value = 3.25
$ws.Range("S5").Value = 1.39
$ws.Range("T5").Value = 2.77
$ws.Range("U5").Value = 2.15
$ws.Range("V5").Value = 1.62
$ws.Range("Y5").Value = 25
$ws.Range("AD5").Value = 8.5
$ws.Range("AE5").Value = 23
$ws.Range("AH5").Value = 5.9
$ws.Range("AI5").Value = 8.5
$ws.Range("AJ5").Value = 8.5
$ws.Range("AO5").Value = 50
$ws.Range("AQ5").Value = 400
$ws.Range("AT5").Value = 2.77
$ws.Range("AX5").Value = 6.4
$ws.Range("AZ5").Value = 19

# Row 7
$ws.Range("G7").Value = 5.5
$ws.Range("H7").Value = 3.65
$ws.Range("I7").Value = 1.55
$ws.Range("J7").Value = 5.7
$ws.Range("K7").Value = 2.15
$ws.Range("L7").Value = 2.12
$ws.Range("N7").Value = 7.3
$ws.Range("P7").Value = 3.2
$ws.Range("Q7").Value = 1.91
$ws.Range("R7").Value = 1.83
$ws.Range("T7").Value = 2.7
$ws.Range("U7").Value = 1.93
$ws.Range("V7").Value = 1.78
$ws.Range("W7").Value = 14.5
$ws.Range("X7").Value = 35
$ws.Range("Y7").Value = 17.5
$ws.Range("Z7").Value = 120
$ws.Range("AA7").Value = 60
$ws.Range("AC7").Value = 7.3
$ws.Range("AG7").Value = 6.2
$ws.Range("AH7").Value = 6.9
$ws.Range("AI7").Value = 8
$ws.Range("AJ7").Value = 11
$ws.Range("AL7").Value = 29
$ws.Range("AN7").Value = 7.1
$ws.Range("AT7").Value = 2.7
$ws.Range("AX7").Value = 7.6
